$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" data row (original row 26) entirely; all following
# rows shift up by one.
$ws.Rows(26).Delete()

# Remove the "SC 92" data row (originally row 28, now row 27 after the
# previous deletion) entirely; all following rows shift up by one more.
$ws.Rows(27).Delete()

# Update imputed/corrected values in column D (and a couple of column C
# cells) to their final values for this seed, matching the new row
# positions after the two row deletions above.
$ws.Range("D2").Value = -13.5
$ws.Range("D6").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("D14").Value = ""
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("D22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("C30").Value = 11.4
$ws.Range("D31").Value = -13.7
$ws.Range("C32").Value = ""
$ws.Range("D33").Value = -14.1
